$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 293.25
$ws.Cells.Item(12, 9).Value = 293.66666
$ws.Cells.Item(12, 10).Value = 292
$ws.Cells.Item(12, 11).Value = 293.66666
$ws.Cells.Item(12, 12).Value = 292
$ws.Cells.Item(12, 13).Value = -123.66666
$ws.Cells.Item(12, 14).Value = -632
$ws.Cells.Item(33, 8).Value = 153.28572
$ws.Cells.Item(33, 9).Value = 154.25
$ws.Cells.Item(33, 11).Value = 154.25
$ws.Cells.Item(33, 13).Value = 74.75
$ws.Cells.Item(51, 8).Value = 3899.8
$ws.Cells.Item(51, 9).Value = 3500
$ws.Cells.Item(51, 10).Value = 4166.3335
$ws.Cells.Item(51, 11).Value = 3500
$ws.Cells.Item(51, 12).Value = 4166.3335
$ws.Cells.Item(51, 13).Value = -3016
$ws.Cells.Item(51, 14).Value = -5134.3335
$ws.Cells.Item(88, 8).Value = 1313.909
$ws.Cells.Item(88, 9).Value = 1946.3334
$ws.Cells.Item(88, 10).Value = 1076.75
$ws.Cells.Item(88, 11).Value = 1946.3334
$ws.Cells.Item(88, 12).Value = 1076.75
$ws.Cells.Item(88, 13).Value = -1540.3334
$ws.Cells.Item(88, 14).Value = -1888.75
$ws.Cells.Item(91, 8).Value = 1313.909
$ws.Cells.Item(91, 9).Value = 1946.3334
$ws.Cells.Item(91, 10).Value = 1076.75
$ws.Cells.Item(91, 11).Value = 1946.3334
$ws.Cells.Item(91, 12).Value = 1076.75
$ws.Cells.Item(91, 13).Value = -542.3334
$ws.Cells.Item(91, 14).Value = -3884.75
$ws.Cells.Item(100, 8).Value = 814.2857
$ws.Cells.Item(100, 9).Value = 750
$ws.Cells.Item(100, 10).Value = 1200
$ws.Cells.Item(100, 11).Value = 750
$ws.Cells.Item(100, 12).Value = 1200
$ws.Cells.Item(100, 13).Value = -209
$ws.Cells.Item(100, 14).Value = -2282
$ws.Cells.Item(135, 8).Value = 4289.154
$ws.Cells.Item(135, 9).Value = 3476.9167
$ws.Cells.Item(135, 11).Value = 31292.2503
$ws.Cells.Item(135, 13).Value = -28757.2503
$ws.Cells.Item(138, 8).Value = 2930.4
$ws.Cells.Item(138, 9).Value = 2163
$ws.Cells.Item(138, 10).Value = 6000
$ws.Cells.Item(138, 11).Value = 6489
$ws.Cells.Item(138, 12).Value = 18000
$ws.Cells.Item(138, 13).Value = -1349
$ws.Cells.Item(138, 14).Value = -28280
$ws.Cells.Item(141, 8).Value = 3664.6667
$ws.Cells.Item(141, 9).Value = 2997.5
$ws.Cells.Item(141, 11).Value = 8992.5
$ws.Cells.Item(141, 13).Value = -3812.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 399.8889
$ws.Cells.Item(4, 9).Value = 268.2
$ws.Cells.Item(4, 10).Value = 564.5
$ws.Cells.Item(4, 11).Value = 268.2
$ws.Cells.Item(4, 12).Value = 564.5
$ws.Cells.Item(4, 13).Value = -152.2
$ws.Cells.Item(4, 14).Value = -796.5
$ws.Cells.Item(32, 8).Value = 3211235.5
$ws.Cells.Item(32, 9).Value = 3184920.5
$ws.Cells.Item(32, 11).Value = 3184920.5
$ws.Cells.Item(32, 13).Value = -3184633.5
$ws.Cells.Item(46, 8).Value = 4805.143
$ws.Cells.Item(46, 9).Value = 4764
$ws.Cells.Item(46, 10).Value = 4812
$ws.Cells.Item(46, 11).Value = 4764
$ws.Cells.Item(46, 12).Value = 4812
$ws.Cells.Item(46, 13).Value = -4445
$ws.Cells.Item(46, 14).Value = -5450
$ws.Cells.Item(61, 8).Value = 2353.5
$ws.Cells.Item(61, 9).Value = 2067.5
$ws.Cells.Item(61, 10).Value = 3497.5
$ws.Cells.Item(61, 11).Value = 2067.5
$ws.Cells.Item(61, 12).Value = 3497.5
$ws.Cells.Item(61, 13).Value = -1855.5
$ws.Cells.Item(61, 14).Value = -3921.5
$ws.Cells.Item(132, 8).Value = 1430.6666
$ws.Cells.Item(132, 9).Value = 1446.5
$ws.Cells.Item(132, 10).Value = 1399
$ws.Cells.Item(132, 11).Value = 4339.5
$ws.Cells.Item(132, 12).Value = 4197
$ws.Cells.Item(132, 13).Value = -1809.5
$ws.Cells.Item(132, 14).Value = -9257
$ws.Cells.Item(136, 8).Value = 2353.5
$ws.Cells.Item(136, 9).Value = 2067.5
$ws.Cells.Item(136, 10).Value = 3497.5
$ws.Cells.Item(136, 11).Value = 6202.5
$ws.Cells.Item(136, 12).Value = 10492.5
$ws.Cells.Item(136, 13).Value = -3652.5
$ws.Cells.Item(136, 14).Value = -15592.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1927.2
$ws.Cells.Item(86, 10).Value = 2074.75
$ws.Cells.Item(86, 12).Value = 2074.75
$ws.Cells.Item(86, 14).Value = -4320.75
$ws.Cells.Item(89, 8).Value = 1927.2
$ws.Cells.Item(89, 10).Value = 2074.75
$ws.Cells.Item(89, 12).Value = 10373.75
$ws.Cells.Item(89, 14).Value = -21605.75
$ws.Cells.Item(105, 8).Value = 2532.6667
$ws.Cells.Item(105, 9).Value = 2532.6667
$ws.Cells.Item(105, 11).Value = 2532.6667
$ws.Cells.Item(105, 13).Value = -785.6667000000002
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 13).ClearContents()
$ws.Cells.Item(134, 8).Value = 1463
$ws.Cells.Item(134, 9).Value = 1447
$ws.Cells.Item(134, 11).Value = 4341
$ws.Cells.Item(134, 13).Value = -1806

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2759.2
$ws.Cells.Item(132, 9).Value = 2449.75
$ws.Cells.Item(132, 11).Value = 7349.25
$ws.Cells.Item(132, 13).Value = -4819.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 890.5
$ws.Cells.Item(122, 10).Value = 890.5
$ws.Cells.Item(122, 12).Value = 8014.5
$ws.Cells.Item(122, 14).Value = -12914.5
$ws.Cells.Item(131, 8).Value = 626834.8
$ws.Cells.Item(131, 9).Value = 1273.2858
$ws.Cells.Item(131, 11).Value = 3819.8574
$ws.Cells.Item(131, 13).Value = 1220.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(98, 8).Value = 15650
$ws.Cells.Item(98, 10).Value = 15650
$ws.Cells.Item(98, 12).Value = 15650
$ws.Cells.Item(98, 14).Value = -21640
$ws.Cells.Item(113, 8).Value = 1000
$ws.Cells.Item(113, 10).Value = 1000
$ws.Cells.Item(113, 12).Value = 1000
$ws.Cells.Item(113, 14).Value = -5340

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2794.5
$ws.Cells.Item(61, 10).Value = 2794.5
$ws.Cells.Item(61, 12).Value = 2794.5
$ws.Cells.Item(61, 14).Value = -3198.5
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).ClearContents()
$ws.Cells.Item(93, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 2794.5
$ws.Cells.Item(113, 10).Value = 2794.5
$ws.Cells.Item(113, 12).Value = 2794.5
$ws.Cells.Item(113, 14).Value = -7134.5
$ws.Cells.Item(136, 8).Value = 3470
$ws.Cells.Item(136, 9).Value = 3702.5
$ws.Cells.Item(136, 10).Value = 3005
$ws.Cells.Item(136, 11).Value = 11107.5
$ws.Cells.Item(136, 12).Value = 9015
$ws.Cells.Item(136, 13).Value = -8557.5
$ws.Cells.Item(136, 14).Value = -14115

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 41667.332
$ws.Cells.Item(62, 10).Value = 55000
$ws.Cells.Item(62, 12).Value = 55000
$ws.Cells.Item(62, 14).Value = -56248
$ws.Cells.Item(65, 8).Value = 41667.332
$ws.Cells.Item(65, 10).Value = 55000
$ws.Cells.Item(65, 12).Value = 275000
$ws.Cells.Item(65, 14).Value = -281240
$ws.Cells.Item(100, 8).Value = 10000238
$ws.Cells.Item(100, 9).Value = 12500242
$ws.Cells.Item(100, 10).Value = 221
$ws.Cells.Item(100, 11).Value = 25000484
$ws.Cells.Item(100, 12).Value = 442
$ws.Cells.Item(100, 13).Value = -24999943
$ws.Cells.Item(100, 14).Value = -1524
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 1981
$ws.Cells.Item(136, 10).Value = 4999
$ws.Cells.Item(136, 12).Value = 14997
$ws.Cells.Item(136, 14).Value = -20097
